$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header cells: "<Name>_old" -> "<Name>_FV2310",
#    "<Name>_new" -> "<Name>_FV2404". Column K ("diff") is left untouched.
# ---------------------------------------------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $base = $cell.Value().ToString() -replace "_old$", ""
    $cell.Value = $base + "_FV2310"
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $base = $cell.Value().ToString() -replace "_new$", ""
    $cell.Value = $base + "_FV2404"
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel table ("Table1") without disturbing
#    the existing header-row style (s=1) - adding a ListObject over a range
#    whose header already carries an explicit style makes Excel capture that
#    style into a new dxf (headerRowDxfId), which the target workbook does
#    not have. So: stash the header formatting on a scratch row, strip the
#    header format, create the table, then restore the formatting from the
#    scratch row and discard it again.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRowIndex = 79
$scratchRange = $ws.Range("A" + $scratchRowIndex + ":U" + $scratchRowIndex)

$ws.Rows(1).Copy() | Out-Null
$ws.Rows($scratchRowIndex).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratchRange.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null                  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows($scratchRowIndex).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
